$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump "Förändrad" (C) date for all data rows 2-36 from 46062 to 46063
$ws.Range("C2:C36").Value = 46063

# Row 12
$ws.Range("A12").Value = "A 17492-2024"
$ws.Range("B12").Value = 45415.50709490741
$ws.Range("G12").Value = 7.7

# Row 13
$ws.Range("A13").Value = "A 28409-2024"
$ws.Range("B13").Value = 45477.62280092593
$ws.Range("G13").Value = 0.4

# Row 14
$ws.Range("A14").Value = "A 46587-2024"
$ws.Range("B14").Value = 45582.76763888889
$ws.Range("G14").Value = 2.8

# Row 15
$ws.Range("A15").Value = "A 46588-2024"
$ws.Range("B15").Value = 45582.77137731481
$ws.Range("G15").Value = 1

# Row 16
$ws.Range("A16").Value = "A 270-2025"
$ws.Range("B16").Value = 45660.48087962963
$ws.Range("G16").Value = 8.9

# Row 17
$ws.Range("A17").Value = "A 28418-2024"
$ws.Range("B17").Value = 45477.62978009259
$ws.Range("G17").Value = 0.2

# Row 18
$ws.Range("A18").Value = "A 50239-2022"
$ws.Range("B18").Value = 44865
$ws.Range("G18").Value = 13.2

# Row 19
$ws.Range("A19").Value = "A 24212-2023"
$ws.Range("B19").Value = 45076
$ws.Range("G19").Value = 5.8

# Row 20
$ws.Range("A20").Value = "A 51434-2025"
$ws.Range("B20").Value = 45949
$ws.Range("G20").Value = 2.8

# Row 21
$ws.Range("A21").Value = "A 52888-2025"
$ws.Range("B21").Value = 45957.56943287037
$ws.Range("G21").Value = 2

# Row 22
$ws.Range("A22").Value = "A 52965-2025"
$ws.Range("B22").Value = 45956
$ws.Range("G22").Value = 0.6

# Row 23
$ws.Range("A23").Value = "A 52960-2025"
$ws.Range("G23").Value = 1.7

# Row 24
$ws.Range("A24").Value = "A 4422-2024"
$ws.Range("B24").Value = 45327.45375
$ws.Range("G24").Value = 4.5

# Row 25
$ws.Range("A25").Value = "A 57001-2025"
$ws.Range("B25").Value = 45977

# Row 26
$ws.Range("A26").Value = "A 59231-2024"
$ws.Range("B26").Value = 45637.58472222222
$ws.Range("G26").Value = 1.3

# Row 28
$ws.Range("A28").Value = "A 49634-2024"
$ws.Range("B28").Value = 45596.59591435185
$ws.Range("G28").Value = 0.7

# Row 29
$ws.Range("A29").Value = "A 49633-2024"
$ws.Range("B29").Value = 45596.59559027778
$ws.Range("G29").Value = 0.8

# Row 30
$ws.Range("A30").Value = "A 46579-2024"
$ws.Range("B30").Value = 45582.75018518518
$ws.Range("G30").Value = 3

# Row 31
$ws.Range("A31").Value = "A 28416-2024"

# Row 32
$ws.Range("A32").Value = "A 21421-2021"
$ws.Range("B32").Value = 44316
$ws.Range("G32").Value = 0.6

# Row 33
$ws.Range("A33").Value = "A 17491-2024"
$ws.Range("B33").Value = 45415.50266203703
$ws.Range("G33").Value = 6.2

# Row 34
$ws.Range("A34").Value = "A 23503-2025"
$ws.Range("B34").Value = 45795
$ws.Range("G34").Value = 14.1

# Row 36
$ws.Range("A36").Value = "A 24086-2025"
$ws.Range("G36").Value = 0.7
